$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 'JAVA (Lab - BCA Lab by MN, SR, VKR, MAN) in 3CM                      '
$ws.Range("C3").Value = 'JAVA (Lab - BCA Lab by MN, SR, VKR, MAN) in 3CM                      '
$ws.Range("B4").Value = 'JAVA (Lab - BCA Lab by KNS, RES, VKR, MAN) in 3CS                    '
$ws.Range("C4").Value = 'JAVA (Lab - BCA Lab by KNS, RES, VKR, MAN) in 3CS                    '
$ws.Range("B5").Value = 'PY (Lab - BCA Lab by KNS, NR, TM, KR, VA) in 3BCA A '
$ws.Range("C5").Value = 'PY (Lab - BCA Lab by KNS, NR, TM, KR, VA) in 3BCA A '
$ws.Range("D5").Value = 'Free                                                '
$ws.Range("E5").Value = 'OOP (Lab - BCA Lab by KNS, LJ) in 1BCA B            '
$ws.Range("F5").Value = 'OOP (Lab - BCA Lab by KNS, LJ) in 1BCA B            '
$ws.Range("G5").Value = 'Free                                                '
$ws.Range("H5").Value = 'Free                                                '
$ws.Range("I5").Value = 'OOP (Lab - BCA Lab by CHA, RV, SD, SH, SA) in 1BCA A'
$ws.Range("J5").Value = 'OOP (Lab - BCA Lab by CHA, RV, SD, SH, SA) in 1BCA A'
$ws.Range("B6").Value = 'PY (Lab - BCA Lab by HU, RA, APR, RES, KR) in 3BCA B '
$ws.Range("C6").Value = 'PY (Lab - BCA Lab by HU, RA, APR, RES, KR) in 3BCA B '
$ws.Range("D6").Value = 'Free                                                 '
$ws.Range("E6").Value = 'Free                                                 '
$ws.Range("F6").Value = 'Free                                                 '
$ws.Range("G6").Value = 'Free                                                 '
$ws.Range("H6").Value = 'Free                                                 '
$ws.Range("I6").Value = 'PSD (Lab - BCA Lab by CHA, FHP, NR, SD, SH) in 1BCA A'
$ws.Range("J6").Value = 'PSD (Lab - BCA Lab by CHA, FHP, NR, SD, SH) in 1BCA A'
$ws.Range("B7").Value = 'MA (Lab - BCA Lab by NR, CYN, SD, SH, MAN) in 3BCA A'
$ws.Range("C7").Value = 'MA (Lab - BCA Lab by NR, CYN, SD, SH, MAN) in 3BCA A'
$ws.Range("D7").Value = 'Free                                                '
$ws.Range("E7").Value = 'MA (Lab - BCA Lab by VKR) in 5BCA A                 '
$ws.Range("F7").Value = 'MA (Lab - BCA Lab by VKR) in 5BCA A                 '
$ws.Range("G7").Value = 'Free                                                '
$ws.Range("H7").Value = 'Free                                                '
$ws.Range("I7").Value = 'WT (Lab - BCA Lab by CHA, FHP, SG, HU, TM) in 1BCA B'
$ws.Range("J7").Value = 'WT (Lab - BCA Lab by CHA, FHP, SG, HU, TM) in 1BCA B'
$ws.Range("B8").Value = 'Free                                                 '
$ws.Range("C8").Value = 'Free                                                 '
$ws.Range("D8").Value = 'Free                                                 '
$ws.Range("E8").Value = 'MA (Lab - BCA Lab by VKR, CYN, SG, HU, RES) in 5BCA A'
$ws.Range("F8").Value = 'MA (Lab - BCA Lab by VKR, CYN, SG, HU, RES) in 5BCA A'
$ws.Range("G8").Value = 'Free                                                 '
$ws.Range("H8").Value = 'Free                                                 '
$ws.Range("I8").Value = 'Free                                                 '
$ws.Range("J8").Value = 'Free                                                 '
$ws.Range("B16").Value = 'PRJ (Lab - MCA Lab by BE, RV, SG, KR, VA) in 5BCA A                    '
$ws.Range("C16").Value = 'PRJ (Lab - MCA Lab by BE, RV, SG, KR, VA) in 5BCA A                    '
$ws.Range("I16").Value = 'PSD (Lab - MCA Lab by AMR, SK, NEB, SME, APR) in 1BCA B                '
$ws.Range("J16").Value = 'PSD (Lab - MCA Lab by AMR, SK, NEB, SME, APR) in 1BCA B                '
$ws.Range("B17").Value = 'WAD (Lab - MCA Lab by LJ, BE, NEB) in 3CS           '
$ws.Range("C17").Value = 'WAD (Lab - MCA Lab by LJ, BE, NEB) in 3CS           '
$ws.Range("D17").Value = 'Free                                                '
$ws.Range("E17").Value = 'MA (Lab - MCA Lab by VKR, CYN, SD, TM, SH) in 3BCA B'
$ws.Range("F17").Value = 'MA (Lab - MCA Lab by VKR, CYN, SD, TM, SH) in 3BCA B'
$ws.Range("G17").Value = 'Free                                                '
$ws.Range("H17").Value = 'Free                                                '
$ws.Range("I17").Value = 'DAS (Lab - MCA Lab by SK, SME, RM, FHP, RA) in 1CM  '
$ws.Range("J17").Value = 'DAS (Lab - MCA Lab by SK, SME, RM, FHP, RA) in 1CM  '
$ws.Range("B18").Value = 'DAP (Lab - MCA Lab by SME, SK, NEB, RM, SR) in 3CM  '
$ws.Range("C18").Value = 'DAP (Lab - MCA Lab by SME, SK, NEB, RM, SR) in 3CM  '
$ws.Range("D18").Value = 'Free                                                '
$ws.Range("E18").Value = 'Free                                                '
$ws.Range("F18").Value = 'Free                                                '
$ws.Range("G18").Value = 'Free                                                '
$ws.Range("H18").Value = 'Free                                                '
$ws.Range("I18").Value = 'OOP (Lab - MCA Lab by KNS, LJ, RV, SG, SA) in 1BCA B'
$ws.Range("J18").Value = 'OOP (Lab - MCA Lab by KNS, LJ, RV, SG, SA) in 1BCA B'
$ws.Range("B19").Value = 'Free                                             '
$ws.Range("C19").Value = 'Free                                             '
$ws.Range("D19").Value = 'Free                                             '
$ws.Range("E19").Value = 'MA (Lab - MCA Lab by MAN) in 5BCA B              '
$ws.Range("F19").Value = 'MA (Lab - MCA Lab by MAN) in 5BCA B              '
$ws.Range("G19").Value = 'Free                                             '
$ws.Range("H19").Value = 'Free                                             '
$ws.Range("I19").Value = 'DCF (Lab - MCA Lab by AMR, BE, MN, RV, SR) in 1CM'
$ws.Range("J19").Value = 'DCF (Lab - MCA Lab by AMR, BE, MN, RV, SR) in 1CM'
$ws.Range("B20").Value = 'Free                                                '
$ws.Range("C20").Value = 'Free                                                '
$ws.Range("D20").Value = 'Free                                                '
$ws.Range("E20").Value = 'MA (Lab - MCA Lab by BE, KNS, RV, NR, MAN) in 5BCA B'
$ws.Range("F20").Value = 'MA (Lab - MCA Lab by BE, KNS, RV, NR, MAN) in 5BCA B'
$ws.Range("G20").Value = 'Free                                                '
$ws.Range("H20").Value = 'Free                                                '
$ws.Range("I20").Value = 'Free                                                '
$ws.Range("J20").Value = 'Free                                                '
$ws.Range("B27").Value = 'PRJ (Lab - BSc Lab by SH) in 5BCA B                 '
$ws.Range("C27").Value = 'PRJ (Lab - BSc Lab by SH) in 5BCA B                 '
$ws.Range("D27").Value = 'Free                                                '
$ws.Range("E27").Value = 'DAS (Lab - BSc Lab by BE, FHP, SD, RES, SH) in 1CS  '
$ws.Range("F27").Value = 'DAS (Lab - BSc Lab by BE, FHP, SD, RES, SH) in 1CS  '
$ws.Range("G27").Value = 'Free                                                '
$ws.Range("H27").Value = 'Free                                                '
$ws.Range("I27").Value = 'WT (Lab - BSc Lab by AMR, BE, LJ, APR, KR) in 1BCA A'
$ws.Range("J27").Value = 'WT (Lab - BSc Lab by AMR, BE, LJ, APR, KR) in 1BCA A'
$ws.Range("B28").Value = 'PRJ (Lab - BSc Lab by RM, LJ, RA, HU, SH) in 5BCA B'
$ws.Range("C28").Value = 'PRJ (Lab - BSc Lab by RM, LJ, RA, HU, SH) in 5BCA B'
$ws.Range("D28").Value = 'Free                                               '
$ws.Range("E28").Value = 'Free                                               '
$ws.Range("F28").Value = 'Free                                               '
$ws.Range("G28").Value = 'Free                                               '
$ws.Range("H28").Value = 'Free                                               '
$ws.Range("I28").Value = 'Free                                               '
$ws.Range("J28").Value = 'Free                                               '
$ws.Range("B29").Value = 'Free                                              '
$ws.Range("C29").Value = 'Free                                              '
$ws.Range("D29").Value = 'Free                                              '
$ws.Range("E29").Value = 'Free                                              '
$ws.Range("F29").Value = 'DA (Lab - BSc Lab by AMR, MN, RV, APR, SA) in 5CME'
$ws.Range("G29").Value = 'DA (Lab - BSc Lab by AMR, MN, RV, APR, SA) in 5CME'
$ws.Range("H29").Value = 'Free                                              '
$ws.Range("I29").Value = 'Free                                              '
$ws.Range("J29").Value = 'Free                                              '
$ws.Range("B32").Value = 'Free                                               '
$ws.Range("C32").Value = 'Free                                               '
$ws.Range("D32").Value = 'WAD (Lab - BSc Lab by CHA, RM, FHP, APR, TM) in 1CM'
$ws.Range("E32").Value = 'WAD (Lab - BSc Lab by CHA, RM, FHP, APR, TM) in 1CM'
$ws.Range("F32").Value = 'DCF (Lab - BSc Lab by RM, LJ, SR, TM, KR) in 1CS   '
$ws.Range("G32").Value = 'DCF (Lab - BSc Lab by RM, LJ, SR, TM, KR) in 1CS   '
$ws.Range("H32").Value = 'Free                                               '
$ws.Range("I32").Value = 'Free                                               '
$ws.Range("J32").Value = 'Free                                               '
